$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = -0.1755
$ws.Range("E2").Value = -0.221
$ws.Range("G2").Value = 0.333084762614964
$ws.Range("H2").Value = 0.333084762614964
$ws.Range("I2").Value = 0.1988565746955009
$ws.Range("J2").Value = 0.170519512801392
$ws.Range("K2").Value = -2.917
$ws.Range("L2").Value = -0.7250807854834702
$ws.Range("M2").Value = 0.07
$ws.Range("N2").Value = 0.00748663101604278
$ws.Range("O2").Value = -0.02399725745629071
$ws.Range("P2").Value = 0.07
$ws.Range("Q2").Value = 0.00748663101604278
$ws.Range("R2").Value = -0.02399725745629071
$ws.Range("U2").Value = 1.759
$ws.Range("V2").Value = 0.1881283422459893
$ws.Range("W2").Value = -0.178324764353042
$ws.Range("X2").Value = 0.1333708057375701
$ws.Range("Y2").Value = -0.3116955700906121
$ws.Range("Z2").Value = 0.1408071121066816
$ws.Range("AA2").Value = 0.02983206425367686
$ws.Range("AB2").Value = 0.0849130204855691
$ws.Range("AC2").Value = -0.05508095623189224
$ws.Range("AD2").Value = 19.38
$ws.Range("AF2").Value = 19.38
$ws.Range("AG2").Value = 17.621
$ws.Range("AH2").Value = 0.6745562130177515
$ws.Range("AI2").Value = 0.5899543378995433
$ws.Range("AJ2").Value = 0.6533313559007823
$ws.Range("AK2").Value = 0.5667556527612492
$ws.Range("AL2").Value = 0.624
$ws.Range("AM2").Value = 0.579
$ws.Range("AN2").Value = 16.85217391304348
$ws.Range("AO2").Value = 1.282051282051282
$ws.Range("AP2").Value = 15.32260869565217
$ws.Range("AQ2").Value = 1.381692573402418

# Row 3
$ws.Range("D3").Value = 0.023
$ws.Range("E3").Value = -0.221
$ws.Range("G3").Value = 0.3753501400560225
$ws.Range("H3").Value = 0.3753501400560225
$ws.Range("I3").Value = 0.2240896358543417
$ws.Range("J3").Value = 0.1602240896358544
$ws.Range("K3").Value = 0.143
$ws.Range("L3").Value = 0.04005602240896358
$ws.Range("M3").Value = 0.07
$ws.Range("N3").Value = 0.01515151515151515
$ws.Range("O3").Value = 0.4895104895104896
$ws.Range("P3").Value = 0.07
$ws.Range("Q3").Value = 0.01515151515151515
$ws.Range("R3").Value = 0.4895104895104896
$ws.Range("U3").Value = 0.719
$ws.Range("V3").Value = 0.1556277056277056
$ws.Range("W3").Value = 0.03666666666666667
$ws.Range("X3").Value = 0.09337797364183388
$ws.Range("Y3").Value = -0.05671130697516721
$ws.Range("Z3").Value = 0.3723792635861062
$ws.Range("AA3").Value = 0.05966412850735372
$ws.Range("AB3").Value = 0.07787630175383953
$ws.Range("AC3").Value = -0.01821217324648581
$ws.Range("AD3").Value = 4.38
$ws.Range("AF3").Value = 4.38
$ws.Range("AG3").Value = 3.661
$ws.Range("AH3").Value = 0.4866666666666666
$ws.Range("AI3").Value = 0.4792122538293216
$ws.Range("AJ3").Value = 0.4420963651732882
$ws.Range("AK3").Value = 0.4347464671654198
$ws.Range("AL3").Value = 0.624
$ws.Range("AM3").Value = 0.579
$ws.Range("AN3").Value = 3.808695652173913
$ws.Range("AO3").Value = 1.282051282051282
$ws.Range("AP3").Value = 3.183478260869566
$ws.Range("AQ3").Value = 1.381692573402418

# Row 4
$ws.Range("D4").Value = -0.374
$ws.Range("K4").Value = -3.06
$ws.Range("L4").Value = -6.754966887417218
$ws.Range("U4").Value = 1.04
$ws.Range("V4").Value = 0.2198731501057082
$ws.Range("W4").Value = -0.3933161953727506
$ws.Range("X4").Value = 0.1733636378333063
$ws.Range("Y4").Value = -0.566679833206057
$ws.Range("Z4").Value = 0.0238621997471555
$ws.Range("AB4").Value = 0.09194973921729868
$ws.Range("AC4").Value = -0.09194973921729868
$ws.Range("AD4").Value = 15
$ws.Range("AF4").Value = 15
$ws.Range("AG4").Value = 13.96
$ws.Range("AH4").Value = 0.7602635580334516
$ws.Range("AI4").Value = 0.6326444538169549
$ws.Range("AJ4").Value = 0.7469234884965222
$ws.Range("AK4").Value = 0.615791795324217
$ws.Range("AM4").Value = 0

# Remove AQ4 entirely (was -0, diff shows the cell is deleted from the row)
$ws.Range("AQ4").ClearContents()
